# Outcomes.xlsx edit: "Included two new moments" — regression of probability
# of expert on past test score (either portfolio or test).
#
# Touches three sheets:
#   - "data"      : updates simulated moments (cols D/E/F, rows 5-22) and
#                   appends two new summary rows (23, 24)
#   - "table"     : (no content change in the diff besides what flows from data)
#   - "table_v2"  : relabels row 20 and moves the active-window focus away
#                   from it back to "data"

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("data")
$wsV2   = $wb.Worksheets.Item("table_v2")

# ---------------------------------------------------------------------------
# 1) "data" sheet — refreshed simulation moments (D/E/F, rows 5-22)
# ---------------------------------------------------------------------------

$wsData.Range("D5").Value = 2.464088200842154
$wsData.Range("E5").Value = 2.581198556184769
$wsData.Range("F5").Value = 0.08182810964922846

$wsData.Range("D6").Value = 0.06927044459805856
$wsData.Range("E6").Value = 0.06250998651981354
$wsData.Range("F6").Value = 0.002812108107274512

$wsData.Range("D7").Value = -0.3994801237015926
$wsData.Range("E7").Value = -0.3443212494254112
$wsData.Range("F7").Value = 0.01137968585875859

$wsData.Range("D8").Value = 0.2199161972911622
$wsData.Range("E8").Value = 0.2285554205775261
$wsData.Range("F8").Value = 0.007737233353145483

$wsData.Range("D9").Value = 2.509349917307424
$wsData.Range("E9").Value = 2.592065150976181
$wsData.Range("F9").Value = 0.08266076218085858

$wsData.Range("D10").Value = 0.2905799883187443
$wsData.Range("E10").Value = 0.2635648551732302
$wsData.Range("F10").Value = 0.01053098017579399

$wsData.Range("D11").Value = 2.490999063512804
$wsData.Range("E11").Value = 2.410592958283658
$wsData.Range("F11").Value = 0.07631052784374842

$wsData.Range("D12").Value = 0.3616012939749292
$wsData.Range("E12").Value = 0.2821355535332045
$wsData.Range("F12").Value = 0.01259466476202391

$wsData.Range("D13").Value = 0.3817226041245452
$wsData.Range("E13").Value = 0.529291248397036
$wsData.Range("F13").Value = 0.01947491925114392

$wsData.Range("D14").Value = 0.1826202992317024
$wsData.Range("E14").Value = 0.1488150260894422
$wsData.Range("F14").Value = 0.008497317180154439

$wsData.Range("D15").Value = 0.2224529920915002
$wsData.Range("E15").Value = 0.2389975308689028
$wsData.Range("F15").Value = 0.01970719446538374

$wsData.Range("D16").Value = 0.1164375179339663
$wsData.Range("E16").Value = 0.226189139557443
$wsData.Range("F16").Value = 0.02150634078225261

$wsData.Range("D17").Value = -0.07164023951563281
$wsData.Range("E17").Value = -0.06684177619523311
$wsData.Range("F17").Value = 0.02030145292122476

$wsData.Range("D18").Value = -0.01872787656146262
$wsData.Range("E18").Value = -0.02534243929760681
$wsData.Range("F18").Value = 0.02070464178194394

$wsData.Range("D19").Value = 0.5639777484659059
$wsData.Range("E19").Value = 0.3161279605158629
$wsData.Range("F19").Value = 0.01064494658365078

$wsData.Range("D20").Value = 0.4813890666782157
$wsData.Range("E20").Value = 0.588282594848905
$wsData.Range("F20").Value = 0.02417879833631291

$wsData.Range("D21").Value = 0.2215609137065696
$wsData.Range("E21").Value = 0.3246538201806384
$wsData.Range("F21").Value = 0.02109262303944222

$wsData.Range("D22").Value = -0.05634598840963018
$wsData.Range("E22").Value = 0.04015706062279335
$wsData.Range("F22").Value = 0.007065554505528104

# ---------------------------------------------------------------------------
# 2) "data" sheet — two new trailing rows summarizing the new moments
# ---------------------------------------------------------------------------

$wsData.Range("C23").Value = "Past portfolio and % expert"
$wsData.Range("E23").Value = 0
$wsData.Range("F23").Value = 0

$wsData.Range("C24").Value = "Past test and % expert"
$wsData.Range("E24").Value = 0
$wsData.Range("F24").Value = 0

# ---------------------------------------------------------------------------
# 3) "table_v2" sheet — relabel the control-group portfolio/test row
# ---------------------------------------------------------------------------

$wsV2.Range("B20").Value = "Average of past portfolio-test (control)"

# ---------------------------------------------------------------------------
# 4) View state — "data" becomes the active/selected sheet & cell, matching
#    the updated selection recorded for each sheet in the workbook.
# ---------------------------------------------------------------------------

$wsV2.Activate()
$wsV2.Range("B21").Select()

$wsData.Activate()
$wsData.Range("C25").Select()
